# Apply the "added headline to program.R files, added Mentoring Program
# 'Balu und Du'" edit to the programs table.
#
# Net content changes in the target workbook:
#   - B3 (program_name for taxReform1990): "Tax Reform 1990" -> "Top Tax Reform 1990"
#   - B4 (program_name for taxReform2000): "Tax Reform 2000" -> "Top Tax Reform 2000"
#   - Current selection moved to B3
#   - Page setup: paper size -> A4 (9), orientation -> portrait

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Top Tax Reform 1990"
$ws.Range("B4").Value = "Top Tax Reform 2000"

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("B3").Select()
